$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1) Refresh the time_taken (column F) values on the "data" sheet - this run
#    of the panel query finished at a later time than the previous one.
# ---------------------------------------------------------------------------
$newTimes = @(
    "2021-10-05 14:33:14.724040",
    "2021-10-05 14:33:14.724048",
    "2021-10-05 14:33:14.724051",
    "2021-10-05 14:33:14.724053",
    "2021-10-05 14:33:14.724056",
    "2021-10-05 14:33:14.724059",
    "2021-10-05 14:33:14.724061",
    "2021-10-05 14:33:14.724064",
    "2021-10-05 14:33:14.724066",
    "2021-10-05 14:33:14.724069",
    "2021-10-05 14:33:14.724071",
    "2021-10-05 14:33:14.724074",
    "2021-10-05 14:33:14.724076",
    "2021-10-05 14:33:14.724079",
    "2021-10-05 14:33:14.724081",
    "2021-10-05 14:33:14.724083",
    "2021-10-05 14:33:14.724086",
    "2021-10-05 14:33:14.724089",
    "2021-10-05 14:33:14.724092",
    "2021-10-05 14:33:14.724094",
    "2021-10-05 14:33:14.724096"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value2 = $newTimes[$i]
}

# ---------------------------------------------------------------------------
# 2) Add a new "metadata" worksheet (placed right after "data") describing
#    the panel query that produced the data tab.
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Match the page margins used by the "data" sheet.
$metaSheet.PageSetup.LeftMargin = $dataSheet.PageSetup.LeftMargin
$metaSheet.PageSetup.RightMargin = $dataSheet.PageSetup.RightMargin
$metaSheet.PageSetup.TopMargin = $dataSheet.PageSetup.TopMargin
$metaSheet.PageSetup.BottomMargin = $dataSheet.PageSetup.BottomMargin
$metaSheet.PageSetup.HeaderMargin = $dataSheet.PageSetup.HeaderMargin
$metaSheet.PageSetup.FooterMargin = $dataSheet.PageSetup.FooterMargin

# Header row (B1:G1) - text labels styled like the "data" sheet header.
$metaSheet.Cells.Item(1, 2).Value2 = "data_name"
$metaSheet.Cells.Item(1, 3).Value2 = "data_id"
$metaSheet.Cells.Item(1, 4).Value2 = "data_version"
$metaSheet.Cells.Item(1, 5).Value2 = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value2 = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value2 = "panel_get_request"

$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row (row 2).
$metaSheet.Cells.Item(2, 1).Value2 = 0
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$metaSheet.Cells.Item(2, 2).Value2 = "Auditory Neuropathy"
$metaSheet.Cells.Item(2, 3).Value2 = 3440

# data_version ("1.1") must stay textual rather than become the number 1.1 -
# round-trip it through a temporary Text-formatted cell so the value lands
# as a string without leaving a custom number-format style behind on the
# target cell itself.
$scratch = $dataSheet.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value2 = "1.1"
$scratch.Copy()
$metaSheet.Range("D2").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false

$metaSheet.Cells.Item(2, 5).Value2 = "2021-01-14T22:03:39.837419Z"
$metaSheet.Cells.Item(2, 6).Value2 = "2021-10-05 14:33:14.720808"
$metaSheet.Cells.Item(2, 7).Value2 = "https://panelapp.agha.umccr.org/api/v1/panels/3440/?format=json"

# Keep "data" as the active/selected sheet, as it was before this edit.
$dataSheet.Select()
$dataSheet.Range("A1").Select()
